$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fiscal position rows for Austria (AT), Belgium (BE) and the
# Netherlands (NL) - "wide spread tests" - appended below the existing
# data (rows 12-17).

# --- Austria ---
$ws.Range("A12").Value = "z0bug.fiscalpos_at_1"
$ws.Range("B12").NumberFormat = "General"
$ws.Range("C12").Value = "z0bug.fiscalpos_at"
$ws.Range("D12").Value = "z0bug.tax_22v"
$ws.Range("E12").Value = "z0bug.tax_eu-1-AT-v"

$ws.Range("A13").Value = "z0bug.fiscalpos_at_2"
$ws.Range("B13").NumberFormat = "General"
$ws.Range("C13").Value = "z0bug.fiscalpos_at"
$ws.Range("D13").Value = "z0bug.tax_10v"
$ws.Range("E13").Value = "z0bug.tax_eu-3-AT-v"

# --- Belgium ---
$ws.Range("A14").Value = "z0bug.fiscalpos_be_1"
$ws.Range("B14").NumberFormat = "General"
$ws.Range("C14").Value = "z0bug.fiscalpos_be"
$ws.Range("D14").Value = "z0bug.tax_22v"
$ws.Range("E14").Value = "z0bug.tax_eu-1-BE-v"

$ws.Range("A15").Value = "z0bug.fiscalpos_be_2"
$ws.Range("B15").NumberFormat = "General"
$ws.Range("C15").Value = "z0bug.fiscalpos_be"
$ws.Range("D15").Value = "z0bug.tax_10v"
$ws.Range("E15").Value = "z0bug.tax_eu-3-BE-v"

# --- Netherlands ---
$ws.Range("A16").Value = "z0bug.fiscalpos_nl_1"
$ws.Range("B16").NumberFormat = "General"
$ws.Range("C16").Value = "z0bug.fiscalpos_nl"
$ws.Range("D16").Value = "z0bug.tax_22v"
$ws.Range("E16").Value = "z0bug.tax_eu-1-NL-v"

$ws.Range("A17").Value = "z0bug.fiscalpos_nl_2"
$ws.Range("B17").NumberFormat = "General"
$ws.Range("C17").Value = "z0bug.fiscalpos_nl"
$ws.Range("D17").Value = "z0bug.tax_10v"
$ws.Range("E17").Value = "z0bug.tax_eu-3-NL-v"

# Widen column E a bit to fit the new longer tax codes.
$ws.Columns.Item(5).ColumnWidth = 17.33

# Match the cursor position left behind by the author after typing the
# new rows.
$ws.Range("E18").Select()
